$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5's quantity/unit/notes area (G5:K5) is being cleared out and the
# merged "notes" block is moved one column to the left (was I5:K5, now G5:I5).

# 1. Break the old I5:K5 merge so the cells can be edited independently.
$ws.Range("I5:K5").UnMerge()

# 2. Clear out the old "unit" (G5, was "штук"), "quantity" (H5, was 5) and
#    "notes" (I5) values - the row no longer carries that data.
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

# 3. Recreate the merge one column to the left: G5:I5 instead of I5:K5.
$ws.Range("G5:I5").Merge()

# 4. The row's font grows from 8pt to 12pt.
$ws.Range("A5:K5").Font.Size = 12
